$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "Deposit / Crypto / ETH" transaction row was inserted as row 2,
# pushing the previously existing rows 2-6 down to rows 3-7.
# Rewrite rows 2-7 (bottom-up) with their final values so no
# pre-existing cell gets clobbered before it is read.

$ws.Range("E7").Value = "Deposit"
$ws.Range("N7").Value = "Crypto"
$ws.Range("P7").Value = "ETH"
$ws.Range("T7").Value = 341.28

$ws.Range("E6").Value = "Deposit"
$ws.Range("N6").Value = "Crypto"
$ws.Range("P6").Value = "ETH"
$ws.Range("T6").Value = 596.35320000000002

$ws.Range("E5").Value = "Withdrawal"
$ws.Range("N5").Value = "Wiretransfer"
$ws.Range("P5").Value = "Anywires"
$ws.Range("T5").Value = 1622.46

$ws.Range("E4").Value = "Withdrawal"
$ws.Range("N4").Value = "Crypto"
$ws.Range("P4").Value = "ETH"
$ws.Range("T4").Value = 500.02510000000001

$ws.Range("E3").Value = "Withdrawal"
$ws.Range("N3").Value = "Crypto"
$ws.Range("P3").Value = "ETH"
$ws.Range("T3").Value = 999.98659999999995

$ws.Range("E2").Value = "Deposit"
$ws.Range("N2").Value = "Crypto"
$ws.Range("P2").Value = "ETH"
$ws.Range("T2").Value = 500.97559999999999

# Update sheet view: scroll so column F is the left-most visible column,
# and the selection now spans the extended data range.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("Q2:S13").Select() | Out-Null
